# Insert two new data rows (new rows 36 and 37) above the existing row 36,
# pushing the former rows 36..75 down to 38..77.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36:A37").EntireRow.Insert()

# ---- New row 36 ----
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 45128
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108001
$ws.Range("J36").Value = "Guayaba"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 230
$ws.Range("N36").Value = 5000
$ws.Range("O36").Value = 6000
$ws.Range("P36").Value = 5652
$ws.Range("Q36").Value = "`$/caja 10 kilos"
$ws.Range("R36").Value = "Región de Arica y Parinacota"
$ws.Range("S36").Value = 565
$ws.Range("T36").Value = 10

# ---- New row 37 ----
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 45128
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100108
$ws.Range("H37").Value = "Tropicales y subtropicales"
$ws.Range("I37").Value = 100108001
$ws.Range("J37").Value = "Guayaba"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 270
$ws.Range("N37").Value = 3000
$ws.Range("O37").Value = 4000
$ws.Range("P37").Value = 3556
$ws.Range("Q37").Value = "`$/caja 10 kilos"
$ws.Range("R37").Value = "Región de Arica y Parinacota"
$ws.Range("S37").Value = 356
$ws.Range("T37").Value = 10
